$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 1.85
$ws.Range("H2").Value = 3.1
$ws.Range("I2").Value = 5
$ws.Range("J2").Value = 2.63
$ws.Range("X2").Value = 7
$ws.Range("Z2").Value = 15
$ws.Range("AG2").Value = 9
$ws.Range("AK2").Value = 51
$ws.Range("AW2").Value = 6.5
$ws.Range("AZ2").Value = 126
$ws.Range("BA2").Value = 201
$ws.Range("G3").Value = 1.55
$ws.Range("H3").Value = 4.1
$ws.Range("I3").Value = 6
$ws.Range("K3").Value = 2.2
$ws.Range("L3").Value = 6.5
$ws.Range("M3").Value = 1.06
$ws.Range("N3").Value = 10
$ws.Range("Q3").Value = 2.05
$ws.Range("R3").Value = 1.75
$ws.Range("X3").Value = 6.5
$ws.Range("AG3").Value = 13
$ws.Range("AH3").Value = 29
$ws.Range("AI3").Value = 19
$ws.Range("AK3").Value = 51
$ws.Range("AO3").Value = 8
$ws.Range("AQ3").Value = 26
$ws.Range("AU3").Value = 9.5
$ws.Range("AW3").Value = 7.5
$ws.Range("AZ3").Value = 151
$ws.Range("BB3").Value = 401
$ws.Range("G4").Value = 1.91
$ws.Range("H4").Value = 3.3
$ws.Range("I4").Value = 4.33
$ws.Range("AC4").Value = 7
$ws.Range("BD4").Value = 151
$ws.Range("G5").Value = 1.9
$ws.Range("I5").Value = 4.33
$ws.Range("J5").Value = 2.6
$ws.Range("M5").Value = 1.06
$ws.Range("N5").Value = 10
$ws.Range("Q5").Value = 2.05
$ws.Range("R5").Value = 1.75
$ws.Range("U5").Value = 1.83
$ws.Range("V5").Value = 1.83
$ws.Range("X5").Value = 8.5
$ws.Range("AA5").Value = 17
$ws.Range("AE5").Value = 15
$ws.Range("AJ5").Value = 41
$ws.Range("AK5").Value = 34
$ws.Range("AV5").Value = 51
$ws.Range("Q6").Value = 2.2
$ws.Range("R6").Value = 1.65
$ws.Range("G7").Value = 4.2
$ws.Range("I7").Value = 1.91
$ws.Range("L7").Value = 2.63
$ws.Range("W7").Value = 9
$ws.Range("X7").Value = 21
$ws.Range("AH7").Value = 7.5
$ws.Range("AJ7").Value = 15
$ws.Range("AO7").Value = 29
$ws.Range("G8").Value = 1.8
$ws.Range("H8").Value = 3.6
$ws.Range("J8").Value = 2.4
$ws.Range("U8").Value = 1.91
$ws.Range("V8").Value = 1.91
$ws.Range("AH8").Value = 23
$ws.Range("AU8").Value = 8.5
$ws.Range("AY8").Value = 34
$ws.Range("BB8").Value = 251
$ws.Range("K10").Value = 1.92
$ws.Range("AJ8").Value = 51
$ws.Range("AK8").Value = 41
